# Villa.xlsx - Task4
#
# The "J.1 / Drywall" task (row 46) was listed before its own predecessors
# J.2 (Interior Brick Wall) and J.3 (Exterior Brick Wall), which are listed
# right after it in rows 47-48. The edit moves the J.1/Drywall row down so
# it sits after J.2 and J.3 (into the previously-blank row 49), leaving
# rows 47, 48 and 50 untouched.
#
# It also updates the sheet's current view: zoom to 125% and select C30
# (replacing the old scrolled-to-A19 / selected-I69 view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move A46:E46 ("J.1", "Drywall", "(9, 18, 27)", "J.2, J.3") down into the
# blank row beneath J.2/J.3 - same as selecting the row, cutting it, and
# pasting it a few rows further down in Excel.
$src = $ws.Range("A46:E46")
$dst = $ws.Range("A49:E49")
[void]$src.Cut($dst)

# Update the view state: select C30 and zoom to 125%.
[void]$ws.Activate()
[void]$ws.Range("C30").Select()
$excel.ActiveWindow.Zoom = 125
